$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.633.85'
$ws.Range("E2").Value = '  +3.58%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.698.93'
$ws.Range("E3").Value = '  +2.37%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.55'
$ws.Range("E5").Value = '  +2.82%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3948'
$ws.Range("E7").Value = '  +1.76%  '
$ws.Range("E8").Value = '  +2.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.548'
$ws.Range("E9").Value = '  +8.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '54.98'
$ws.Range("E10").Value = '  +11.83%  '
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08817'
$ws.Range("E12").Value = '  +2.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.300'
$ws.Range("E13").Value = '  +8.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.46'
$ws.Range("E14").Value = '  +3.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001333'
$ws.Range("E15").Value = '  +2.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.666'
$ws.Range("E16").Value = '  +6.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.699.65'
$ws.Range("E17").Value = '  +2.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '101.36'
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07104'
$ws.Range("E19").Value = '  +4.64%  '
$ws.Range("E20").Value = '  +4.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.910'
$ws.Range("E21").Value = '  +4.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.18'
$ws.Range("E23").Value = '  +3.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.618.88'
$ws.Range("E24").Value = '  +3.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.065'
$ws.Range("E25").Value = '  +12.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.336'
$ws.Range("E26").Value = '  +1.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.47'
$ws.Range("E27").Value = '  +3.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.93'
$ws.Range("E28").Value = '  +1.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.252'
$ws.Range("E29").Value = '  +1.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.28'
$ws.Range("E30").Value = '  +3.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.650'
$ws.Range("E31").Value = '  +16.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.885.96'
$ws.Range("E32").Value = '  +2.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.108'
$ws.Range("E33").Value = '  -2.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.448'
$ws.Range("E34").Value = '  +14.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08563'
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.31'
$ws.Range("E36").Value = '  +9.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2760'
$ws.Range("E37").Value = '  +4.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.945'
$ws.Range("E38").Value = '  -1.77%  '
$ws.Range("E39").Value = '  +3.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02806'
$ws.Range("E40").Value = '  +11.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09089'
$ws.Range("E41").Value = '  +3.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7775'
$ws.Range("E42").Value = '  +3.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.465'
$ws.Range("E43").Value = '  +1.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7300'
$ws.Range("E44").Value = '  +4.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.53'
$ws.Range("E45").Value = '  +5.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.522'
$ws.Range("E46").Value = '  +6.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.214'
$ws.Range("E47").Value = '  +4.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.395'
$ws.Range("E48").Value = '  +21.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.001'
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.99'
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08043'
$ws.Range("E51").Value = '  +4.13%  '
